$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws1.Range("A1").Value = "test"
